# Clear the "{}" placeholder values from the query (E) and body (F) columns
# for rows that previously held the shared-string "{}" value, leaving the
# cells blank so Excel can match empty query/body in tests.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellsToClear = @(
    "E2", "F2",
    "E3", "F3",
    "E4",
    "F5",
    "E6", "F6",
    "F7",
    "E8", "F8",
    "E9", "F9",
    "E10",
    "E11",
    "E12", "F12",
    "E13", "F13",
    "E14", "F14"
)

foreach ($addr in $cellsToClear) {
    $ws.Range($addr).ClearContents()
}

$ws.Range("F14").Select()
